$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (data starts at row 2).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

# Every data row's "Förändrad" (column C) date stamp moves from 45192 to 45202.
for ($r = 2; $r -le $lastRow; $r++) {
    $c = $ws.Cells.Item($r, 3)
    if ($c.Value2 -eq 45192) {
        $c.Value2 = 45202
    }
}

# Row 3 ("A 42296-2019") gained a new signal species ("Dropptaggsvamp"),
# bumping the Signalarter (I) and Alla arter (Q) counts by one.
$ws.Cells.Item(3, 9).Value2 = 2
$ws.Cells.Item(3, 17).Value2 = 9

$speciesCell = $ws.Cells.Item(3, 18)
$species = $speciesCell.Value2
$species = $species.Replace("Vedskivlav`r`nSkuggblåslav", "Vedskivlav`r`nDropptaggsvamp`r`nSkuggblåslav")
$speciesCell.Value2 = $species
